$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (r = D, M, N, O, P, Q, S, T), taken from the
# corresponding source row of the previous data (rows are permuted).
$data = @{
    2 = @(44491, 300, 14000, 15000, 14500, "`$/bandeja 10 kilos", 1450, 10)
    3 = @(44418, 240, 10000, 11000, 10500, "`$/bandeja 10 kilos", 1050, 10)
    4 = @(44263, 250, 21000, 22000, 21500, "`$/caja 18 kilos",    1194, 18)
    5 = @(44307, 250, 19000, 20000, 19500, "`$/bandeja 18 kilos", 1083, 18)
    6 = @(44291, 200, 17000, 18000, 17500, "`$/bandeja 18 kilos",  972, 18)
    7 = @(44487, 300, 14000, 15000, 14500, "`$/bandeja 10 kilos", 1450, 10)
    8 = @(44489, 300, 26000, 27000, 26500, "`$/bandeja 18 kilos", 1472, 18)
    9 = @(44323, 270, 21000, 22000, 21500, "`$/bandeja 18 kilos", 1194, 18)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]

    $ws.Range("D$r").Value = $vals[0]
    $ws.Range("M$r").Value = $vals[1]
    $ws.Range("N$r").Value = $vals[2]
    $ws.Range("O$r").Value = $vals[3]
    $ws.Range("P$r").Value = $vals[4]
    $ws.Range("Q$r").Value = $vals[5]
    $ws.Range("S$r").Value = $vals[6]
    $ws.Range("T$r").Value = $vals[7]
}
